$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new date cells (column A) the same number format as the
# preceding date rows (17-20) before assigning values, so they reuse the
# existing date style instead of Excel auto-creating a new one.
$ws.Range("A21:A24").NumberFormat = $ws.Range("A20").NumberFormat

# New rows of time-tracking data for labo 3 entries
$ws.Range("A21").Value = Get-Date -Year 2024 -Month 8 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("B21").Value = 7.5
$ws.Range("C21").Value = "Structure générale labo 3 et code function single_population_growth"

$ws.Range("A22").Value = Get-Date -Year 2024 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "Terminer exercice 1 du labo 3"

$ws.Range("A23").Value = Get-Date -Year 2024 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = "Terminer exercice 2,3,4,5,A,B du labo 3"

$ws.Range("A24").Value = Get-Date -Year 2024 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("B24").Value = 3.5
$ws.Range("C24").Value = "Terminer labo 3 et création guide introduction modélisation avec R"

$ws.Range("C24").Select()

$wb.Save()
